$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.084879506803141
$ws.Range("D2").Value = 1.038810547279591
$ws.Range("E2").Value = 1.085398911711135
$ws.Range("F2").Value = 1.091478575385711
$ws.Range("I2").Value = 1.038008778151634
$ws.Range("J2").Value = 1.089737117956965
$ws.Range("K2").Value = 1.041597520056828
$ws.Range("L2").Value = 1.088059443004813
$ws.Range("M2").Value = 1.094123465823692
$ws.Range("N2").Value = 1.091284669659058

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.088975364921588
$ws.Range("D3").Value = 1.039552306659291
$ws.Range("E3").Value = 1.08903124411268
$ws.Range("F3").Value = 1.095186091022087
$ws.Range("I3").Value = 1.038397919137442
$ws.Range("J3").Value = 1.093481699485741
$ws.Range("K3").Value = 1.042149549683819
$ws.Range("L3").Value = 1.091504218873795
$ws.Range("M3").Value = 1.097644451205149
$ws.Range("N3").Value = 1.095034568923114

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.091595031847798
$ws.Range("D4").Value = 1.040025909496182
$ws.Range("E4").Value = 1.091353050600565
$ws.Range("F4").Value = 1.097553529568642
$ws.Range("I4").Value = 1.038642292731002
$ws.Range("J4").Value = 1.095874777686565
$ws.Range("K4").Value = 1.04249973553488
$ws.Range("L4").Value = 1.09370464723996
$ws.Range("M4").Value = 1.099891189306663
$ws.Range("N4").Value = 1.097431045569474

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.092689243563364
$ws.Range("D5").Value = 1.040223511489554
$ws.Range("E5").Value = 1.092322511795931
$ws.Range("F5").Value = 1.098541459665494
$ws.Range("I5").Value = 1.038743271890845
$ws.Range("J5").Value = 1.09687388509596
$ws.Range("K5").Value = 1.042645297288717
$ws.Range("L5").Value = 1.09462307074679
$ws.Range("M5").Value = 1.100828369573772
$ws.Range("N5").Value = 1.098431571825989

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.092872557105121
$ws.Range("D6").Value = 1.040256602390052
$ws.Range("E6").Value = 1.092484905999227
$ws.Range("F6").Value = 1.098706913167426
$ws.Range("I6").Value = 1.038760124427132
$ws.Range("J6").Value = 1.097041238741462
$ws.Range("K6").Value = 1.042669641262086
$ws.Range("L6").Value = 1.094776894710153
$ws.Range("M6").Value = 1.100985301333921
$ws.Range("N6").Value = 1.098599163132864

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.091609680334241
$ws.Range("D7").Value = 1.040028555729837
$ws.Range("E7").Value = 1.091366030351441
$ws.Range("F7").Value = 1.097566758901082
$ws.Range("I7").Value = 1.038643648888064
$ws.Range("J7").Value = 1.095888154802663
$ws.Range("K7").Value = 1.042501687016725
$ws.Range("L7").Value = 1.093716945071645
$ws.Range("M7").Value = 1.09990374055503
$ws.Range("N7").Value = 1.097444441682611

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.08627020751111
$ws.Range("D8").Value = 1.039062560557335
$ws.Range("E8").Value = 1.086632516861617
$ws.Range("F8").Value = 1.092738218206823
$ws.Range("I8").Value = 1.038141841749576
$ws.Range("J8").Value = 1.09100894580087
$ws.Range("K8").Value = 1.041785547481764
$ws.Range("L8").Value = 1.089229662064851
$ws.Range("M8").Value = 1.095320066392752
$ws.Range("N8").Value = 1.092558303644382

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.076615796294666
$ws.Range("D9").Value = 1.037310509578162
$ws.Range("E9").Value = 1.078062959753311
$ws.Range("F9").Value = 1.083977756746222
$ws.Range("I9").Value = 1.037199653689209
$ws.Range("J9").Value = 1.082171791088011
$ws.Range("K9").Value = 1.040468818223279
$ws.Range("L9").Value = 1.081094163641789
$ws.Range("M9").Value = 1.086991377974146
$ws.Range("N9").Value = 1.083708599158196

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.069998379118075
$ws.Range("D10").Value = 1.036107347683629
$ws.Range("E10").Value = 1.072182018535426
$ws.Range("F10").Value = 1.077953252715062
$ws.Range("I10").Value = 1.036531009732928
$ws.Range("J10").Value = 1.076104505182689
$ws.Range("K10").Value = 1.03955256276124
$ws.Range("L10").Value = 1.075503135914606
$ws.Range("M10").Value = 1.081255362524539
$ws.Range("N10").Value = 1.07763269701096

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067086233371979
$ws.Range("D11").Value = 1.035577648003339
$ws.Range("E11").Value = 1.069592328018854
$ws.Range("F11").Value = 1.075297406915423
$ws.Range("I11").Value = 1.036231496529989
$ws.Range("J11").Value = 1.073432089192601
$ws.Range("K11").Value = 1.039146317994334
$ws.Range("L11").Value = 1.073039209055765
$ws.Range("M11").Value = 1.078724691401314
$ws.Range("N11").Value = 1.074956485883637

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.065997169480866
$ws.Range("D12").Value = 1.035379548455493
$ws.Range("E12").Value = 1.068623609650752
$ws.Range("F12").Value = 1.074303504248322
$ws.Range("I12").Value = 1.036118709882874
$ws.Range("J12").Value = 1.072432323237293
$ws.Range("K12").Value = 1.038993958834565
$ws.Range("L12").Value = 1.072117248223402
$ws.Range("M12").Value = 1.077777332112458
$ws.Range("N12").Value = 1.073955300145997

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066231116198979
$ws.Range("D13").Value = 1.035422102951317
$ws.Range("E13").Value = 1.06883171537914
$ws.Range("F13").Value = 1.07451703994968
$ws.Range("I13").Value = 1.036142972986677
$ws.Range("J13").Value = 1.072647103470495
$ws.Range("K13").Value = 1.039026707113959
$ws.Range("L13").Value = 1.072315322194588
$ws.Range("M13").Value = 1.077980881778464
$ws.Range("N13").Value = 1.074170385391765

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.066996363232777
$ws.Range("D14").Value = 1.035561300689946
$ws.Range("E14").Value = 1.069512393803589
$ws.Range("F14").Value = 1.075215403497487
$ws.Range("I14").Value = 1.036222205047924
$ws.Range("J14").Value = 1.073349595176366
$ws.Range("K14").Value = 1.039133753964179
$ws.Range("L14").Value = 1.072963138896663
$ws.Range("M14").Value = 1.078646534240618
$ws.Range("N14").Value = 1.074873874716435

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067466871088977
$ws.Range("D15").Value = 1.035646885639146
$ws.Range("E15").Value = 1.069930872987618
$ws.Range("F15").Value = 1.075644697710983
$ws.Range("I15").Value = 1.036270818209021
$ws.Range("J15").Value = 1.073781471465241
$ws.Range("K15").Value = 1.039199514300641
$ws.Range("L15").Value = 1.073361376900734
$ws.Range("M15").Value = 1.079055680710644
$ws.Range("N15").Value = 1.075306364319177

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.070190634967806
$ws.Range("D16").Value = 1.036142315384538
$ws.Range("E16").Value = 1.072352952137795
$ws.Range("F16").Value = 1.078128491753353
$ws.Range("I16").Value = 1.036550674198883
$ws.Range("J16").Value = 1.076280884937667
$ws.Range("K16").Value = 1.039579320831304
$ws.Range("L16").Value = 1.075665728575035
$ws.Range("M16").Value = 1.081422300334221
$ws.Range("N16").Value = 1.077809327245421

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.071886427437968
$ws.Range("D17").Value = 1.036450725943071
$ws.Range("E17").Value = 1.073860483426793
$ws.Range("F17").Value = 1.079673658525686
$ws.Range("I17").Value = 1.036723523628469
$ws.Range("J17").Value = 1.077836371237646
$ws.Range("K17").Value = 1.039814996224316
$ws.Range("L17").Value = 1.077099480036738
$ws.Range("M17").Value = 1.082894041979296
$ws.Range("N17").Value = 1.079367022514363

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072871066465231
$ws.Range("D18").Value = 1.03662977746123
$ws.Range("E18").Value = 1.074735653430601
$ws.Range("F18").Value = 1.080570396353827
$ws.Range("I18").Value = 1.036823382656899
$ws.Range("J18").Value = 1.078739317003574
$ws.Range("K18").Value = 1.039951547498354
$ws.Range("L18").Value = 1.077931637010522
$ws.Range("M18").Value = 1.083747976342384
$ws.Range("N18").Value = 1.080271250566849

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073206051676999
$ws.Range("D19").Value = 1.036690688100764
$ws.Range("E19").Value = 1.075033369201388
$ws.Range("F19").Value = 1.080875401349072
$ws.Range("I19").Value = 1.036857269979698
$ws.Range("J19").Value = 1.079046470819426
$ws.Range("K19").Value = 1.039997953906867
$ws.Range("L19").Value = 1.078214689832484
$ws.Range("M19").Value = 1.084038390286606
$ws.Range("N19").Value = 1.08057884057635

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.071704951564205
$ws.Range("D20").Value = 1.036417723452654
$ws.Range("E20").Value = 1.073699170710889
$ws.Range("F20").Value = 1.079508347711822
$ws.Range("I20").Value = 1.036705078172441
$ws.Range("J20").Value = 1.077669933729416
$ws.Range("K20").Value = 1.039789805288013
$ws.Range("L20").Value = 1.076946080911823
$ws.Range("M20").Value = 1.082736606420949
$ws.Range("N20").Value = 1.079200348645781

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066771223116088
$ws.Range("D21").Value = 1.035520347858675
$ws.Range("E21").Value = 1.069312140957718
$ws.Range("F21").Value = 1.075009959794417
$ws.Range("I21").Value = 1.03619891581051
$ws.Range("J21").Value = 1.073142927800019
$ws.Range("K21").Value = 1.039102271973123
$ws.Range("L21").Value = 1.072772561739974
$ws.Range("M21").Value = 1.078450721737165
$ws.Range("N21").Value = 1.074666913848709

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.063626410868246
$ws.Range("D22").Value = 1.034948328802041
$ws.Range("E22").Value = 1.066514384867867
$ws.Range("F22").Value = 1.072138645482657
$ws.Range("I22").Value = 1.03587177747157
$ws.Range("J22").Value = 1.070255306398869
$ws.Range("K22").Value = 1.038661516714318
$ws.Range("L22").Value = 1.070109303518949
$ws.Range("M22").Value = 1.075713298859841
$ws.Range("N22").Value = 1.071775191693955

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.065297707764219
$ws.Range("D23").Value = 1.035252318796615
$ws.Range("E23").Value = 1.06800137272257
$ws.Range("F23").Value = 1.073664967819071
$ws.Range("I23").Value = 1.036046054519033
$ws.Range("J23").Value = 1.071790114286683
$ws.Range("K23").Value = 1.038895984829225
$ws.Range("L23").Value = 1.071524964239248
$ws.Range("M23").Value = 1.077168612850808
$ws.Range("N23").Value = 1.073312179185015

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.071786966545837
$ws.Range("D24").Value = 1.03643263844845
$ws.Range("E24").Value = 1.073772073779601
$ws.Range("F24").Value = 1.079583058535718
$ws.Range("I24").Value = 1.036713415847831
$ws.Range("J24").Value = 1.077745153074621
$ws.Range("K24").Value = 1.039801190810618
$ws.Range("L24").Value = 1.077015408103551
$ws.Range("M24").Value = 1.082807758681778
$ws.Range("N24").Value = 1.079275674811084

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.07914240065816
$ws.Range("D25").Value = 1.037769524430211
$ws.Range("E25").Value = 1.080306903480079
$ws.Range("F25").Value = 1.086273885183779
$ws.Range("I25").Value = 1.037450253768617
$ws.Range("J25").Value = 1.084486266824293
$ws.Range("K25").Value = 1.040815876895453
$ws.Range("L25").Value = 1.083225829821006
$ws.Range("M25").Value = 1.089175803823032
$ws.Range("N25").Value = 1.086026361715497
